# Refactor CO2 in constants: rename the Cyrillic "СО2" label to the Latin "CO2".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Cell B12 currently holds the Cyrillic-lettered "СО2" designation for the
# "Концентрация углекислого газа в атмосфере" row; rename it to the Latin "CO2".
$ws.Range("B12").Value = "CO2"

# Reflect the edit as the active selection, matching the authored change.
$ws.Range("B12").Select()
